# Auto-generated edit script: updates crypto price/volume table cells
# to match the scraped values from Mon May  6 04:43:46 UTC 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163. Used via Formula+Copy+PasteSpecial round-trip
# so numeric-looking text (e.g. "5.90", "10.00") is written as an exact
# text string instead of being auto-coerced to a Double (which would
# silently drop meaningful trailing zeros) -- and crucially this path
# leaves no NumberFormat/style residue behind on the cell.
$xlPasteValues = -4163

$ws.Range('D2').Value = '64.166.92'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '3.142.59'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = '="590.49"'
$ws.Range('D5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Formula = '="146.91"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.133.64'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('E10').Value = '  +1.87%  '
$ws.Range('D11').Formula = '="5.90"'
$ws.Range('D11').Copy() | Out-Null
$ws.Range('D11').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E11').Value = '  +5.26%  '
$ws.Range('D12').Formula = '="0.457"'
$ws.Range('D12').Copy() | Out-Null
$ws.Range('D12').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Formula = '="0.0000247"'
$ws.Range('D13').Copy() | Out-Null
$ws.Range('D13').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Formula = '="37.58"'
$ws.Range('D14').Copy() | Out-Null
$ws.Range('D14').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').Value = '3.660.09'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Formula = '="7.26"'
$ws.Range('D17').Copy() | Out-Null
$ws.Range('D17').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').Value = '63.900.56'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '3.137.41'
$ws.Range('E19').Value = '  +1.30%  '
$ws.Range('D20').Formula = '="469.89"'
$ws.Range('D20').Copy() | Out-Null
$ws.Range('D20').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Formula = '="14.40"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').Formula = '="0.731"'
$ws.Range('D22').Copy() | Out-Null
$ws.Range('D22').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').Formula = '="7.57"'
$ws.Range('D23').Copy() | Out-Null
$ws.Range('D23').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('D24').Formula = '="2.40"'
$ws.Range('D24').Copy() | Out-Null
$ws.Range('D24').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E24').Value = '  +12.90%  '
$ws.Range('D25').Formula = '="13.12"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('D26').Formula = '="80.98"'
$ws.Range('D26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Formula = '="10.00"'
$ws.Range('D28').Copy() | Out-Null
$ws.Range('D28').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E28').Value = '  +11.54%  '
$ws.Range('D29').Formula = '="2.71"'
$ws.Range('D29').Copy() | Out-Null
$ws.Range('D29').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').Value = '  +6.18%  '
$ws.Range('E32').Value = '  +0.79%  '
$ws.Range('D33').Formula = '="0.114"'
$ws.Range('D33').Copy() | Out-Null
$ws.Range('D33').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E33').Value = '  +3.79%  '
$ws.Range('D34').Formula = '="27.64"'
$ws.Range('D34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E34').Value = '  +3.86%  '
$ws.Range('D35').Value = '0.0₃0853'
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  -3.71%  '
$ws.Range('D40').Formula = '="459.56"'
$ws.Range('D40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E40').Value = '  +5.83%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Formula = '="51.25"'
$ws.Range('D41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Formula = '="9.32"'
$ws.Range('D42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E42').Value = '  +6.93%  '
$ws.Range('E43').Value = '  +7.63%  '
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').Value = '2.888.80'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').Formula = '="39.92"'
$ws.Range('D46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E46').Value = '  +11.39%  '
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').Formula = '="133.37"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial($xlPasteValues) | Out-Null
$ws.Range('E48').Value = '  +8.16%  '
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('E51').Value = '  +3.62%  '
